$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (superficie-util-agrupada): dim -> medida ; skos:Concept -> xsd:int ; remove mapping file row
$ws.Range("A2").Value = "iaest-measure:superficie-util-agrupada"
$ws.Range("A3").Value = "medida"
$ws.Range("A4").Value = "xsd:int"
$ws.Range("A5").Clear()

# Column D (municipio-nombre): measure -> dimension ; medida -> dim ; xsd:int -> URI-Municipio
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"

# Column I (superficie-util): measure -> dimension ; medida -> dim ; xsd:int -> skos:Concept ; add mapping file
$ws.Range("I2").Value = "iaest-dimension:superficie-util"
$ws.Range("I3").Value = "dim"
$ws.Range("I4").Value = "skos:Concept"
$ws.Range("B5").Copy($ws.Range("I5"))
$ws.Range("I5").Value = "mapping-superficie-util.xlsx"
